# Reverse the order of comma-separated names in the "Recorded By" column (G)
# for every data row on the active worksheet. Cells containing a single
# name (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(",")

        if ($parts.Length -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            $reversed = $trimmed[($trimmed.Length - 1)..0]

            $newVal = [string]::Join(", ", $reversed)
            $cell.Value = $newVal
        }
    }
}
